$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '44.218.10'
$ws.Range('E2').Value = '  +4.69%  '

$ws.Range('D3').Value = '2.224.85'
$ws.Range('E3').Value = '  +2.51%  '

$ws.Range('E4').Value = '  +0.09%  '

$ws.Range('D5').Value = "'260.37"
$ws.Range('E5').Value = '  +2.91%  '

$ws.Range('D6').Value = "'83.18"
$ws.Range('E6').Value = '  +13.23%  '

$ws.Range('D7').Value = "'0.630"
$ws.Range('E7').Value = '  +3.31%  '

$ws.Range('D9').Value = "'0.611"
$ws.Range('E9').Value = '  +5.09%  '

$ws.Range('D10').Value = "'44.34"
$ws.Range('E10').Value = '  +10.05%  '

$ws.Range('D11').Value = "'0.0934"
$ws.Range('E11').Value = '  +2.73%  '

$ws.Range('E12').Value = '  +4.58%  '

$ws.Range('D13').Value = "'0.104"
$ws.Range('E13').Value = '  +3.40%  '

$ws.Range('D14').Value = '2.561.18'
$ws.Range('E14').Value = '  +2.51%  '

$ws.Range('D15').Value = "'14.67"
$ws.Range('E15').Value = '  +3.32%  '

$ws.Range('D16').Value = '2.223.20'
$ws.Range('E16').Value = '  +1.17%  '

$ws.Range('D17').Value = "'0.781"
$ws.Range('E17').Value = '  +2.47%  '

$ws.Range('D18').Value = '44.093.99'
$ws.Range('E18').Value = '  +4.73%  '

$ws.Range('D19').Value = "'0.0000104"
$ws.Range('E19').Value = '  +2.00%  '

$ws.Range('D20').Value = "'71.75"
$ws.Range('E20').Value = '  +1.64%  '

$ws.Range('E21').Value = '  +3.34%  '

$ws.Range('E22').Value = '  +9.92%  '

$ws.Range('D23').Value = "'233.88"
$ws.Range('E23').Value = '  +3.32%  '

$ws.Range('D24').Value = "'9.28"
$ws.Range('E24').Value = '  -2.72%  '

$ws.Range('E25').Value = '  +0.09%  '

$ws.Range('D26').Value = "'10.83"
$ws.Range('E26').Value = '  +3.68%  '

$ws.Range('D27').Value = "'41.07"
$ws.Range('E27').Value = '  +12.05%  '

$ws.Range('E28').Value = '  +1.44%  '

$ws.Range('E29').Value = '  +2.61%  '

$ws.Range('D30').Value = "'2.23"
$ws.Range('E30').Value = '  +3.99%  '

$ws.Range('E31').Value = '  +2.65%  '

$ws.Range('B32').Value = 'EthereumClassic'
$ws.Range('C32').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D32').Value = "'20.69"
$ws.Range('E32').Value = '  +3.46%  '

$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').Value = "'0.0888"
$ws.Range('E33').Value = '  +10.44%  '

$ws.Range('E34').Value = '  +4.64%  '

$ws.Range('E35').Value = '  +9.53%  '

$ws.Range('E36').Value = '  +2.34%  '

$ws.Range('D37').Value = "'0.0364"
$ws.Range('E37').Value = '  +10.93%  '

$ws.Range('E38').Value = '  +7.50%  '

$ws.Range('E39').Value = '  +12.12%  '

$ws.Range('D40').Value = "'2.96"
$ws.Range('E40').Value = '  +22.10%  '

$ws.Range('D41').Value = "'2.12"
$ws.Range('E41').Value = '  +3.99%  '

$ws.Range('D42').Value = "'63.86"
$ws.Range('E42').Value = '  +8.33%  '

$ws.Range('D43').Value = "'5.57"
$ws.Range('E43').Value = '  +9.17%  '

$ws.Range('E44').Value = '  +3.35%  '

$ws.Range('D45').Value = "'103.20"
$ws.Range('E45').Value = '  +0.95%  '

$ws.Range('D46').Value = "'0.0991"
$ws.Range('E46').Value = '  +2.37%  '

$ws.Range('D47').Value = "'8.37"
$ws.Range('E47').Value = '  +1.17%  '

$ws.Range('B48').Value = 'Stacks'
$ws.Range('C48').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D48').Value = "'1.57"
$ws.Range('E48').Value = '  +29.22%  '

$ws.Range('B49').Value = 'ARBITRUM'
$ws.Range('C49').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D49').Value = "'1.13"
$ws.Range('E49').Value = '  +3.69%  '

$ws.Range('B50').Value = 'TrustWalletToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D50').Value = "'1.17"
$ws.Range('E50').Value = '  +3.96%  '

$ws.Range('B51').Value = 'WOONetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range('D51').Value = "'0.443"
$ws.Range('E51').Value = '  -5.08%  '

# Reset number format/style on cells forced to text, so style index stays default (0)
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D7').Style = 'Normal'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D11').Style = 'Normal'
$ws.Range('D13').Style = 'Normal'
$ws.Range('D15').Style = 'Normal'
$ws.Range('D17').Style = 'Normal'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D23').Style = 'Normal'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D50').Style = 'Normal'
$ws.Range('D51').Style = 'Normal'
